# Updated symbol list on Wed Jan 25 19:32:05 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2"  = "301.63"
    "E2"  = "-2.80%"
    "E3"  = "-0.26%"
    "E4"  = "-0.59%"
    "D5"  = "0.08002"
    "E5"  = "-2.66%"
    "D6"  = "1.930"
    "E6"  = "-6.34%"
    "D7"  = "7.768"
    "E7"  = "-2.08%"
    "D8"  = "0.9276"
    "E8"  = "0.04%"
    "D9"  = "0.1557"
    "E9"  = "37.56%"
    "D10" = "0.1897"
    "E10" = "-1.19%"
    "D11" = "0.09016"
    "E11" = "-3.13%"
    "D12" = "0.03431"
    "E12" = "-4.09%"
    "D13" = "0.09890"
    "E13" = "-0.18%"
    "D14" = "0.001402"
    "E14" = "-2.61%"
    "D15" = "0.005750"
    "E15" = "0.21%"
    "D16" = "3.535"
    "E16" = "1.95%"
    "D17" = "4.042"
    "D18" = "2.965"
    "E18" = "-0.70%"
    "D19" = "0.3445"
    "E19" = "0.48%"
    "E20" = "-0.51%"
    "D21" = "5.029"
    "E21" = "-1.35%"
    "D23" = "0.04493"
    "E23" = "-1.00%"
    "D24" = "0.001211"
    "E24" = "-1.23%"
    "D25" = "0.004771"
    "E25" = "-0.88%"
    "D26" = "0.0001230"
    "E26" = "-1.68%"
    "D27" = "0.0003022"
    "E27" = "-32.02%"
    "D39" = "0.01846"
    "E39" = "-7.20%"
    "D40" = "0.04762"
    "E40" = "-3.58%"
    "D41" = "0.01060"
    "E41" = "6.94%"
    "D42" = "0.007351"
    "E42" = "-3.84%"
    "E43" = "-4.13%"
    "D44" = "0.002109"
    "E44" = "-0.63%"
    "D45" = "0.009714"
    "E45" = "-16.17%"
    "D46" = "0.00006231"
    "E46" = "-4.94%"
    "D47" = "0.00000000750"
    "E47" = "-0.09%"
    "E48" = "-64.06%"
    "D50" = "0.00002099"
    "E50" = "-0.09%"
    "D51" = "0.0001999"
    "E51" = "-0.09%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
